# The author removed the post row for "「開かれたターと結ばれたター」"
# (originally row 533) from the sheet, which shifts every following
# row up by one (534->533, 535->534, ... 646->645) and shrinks the
# sheet's used range from A1:C646 to A1:C645.
#
# Deleting the whole worksheet row reproduces exactly that: Excel
# removes the row and renumbers/shifts everything below it up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(533).Delete()
